$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2632700.2
$ws.Range("I19").Value = 3760154.2
$ws.Range("J19").Value = 1974
$ws.Range("K19").Value = 3760154.2
$ws.Range("L19").Value = 1974
$ws.Range("M19").Value = -3759979.2
$ws.Range("N19").Value = -2324

$ws.Range("H33").Value = 330.52942
$ws.Range("I33").Value = 309.23077
$ws.Range("K33").Value = 309.23077
$ws.Range("M33").Value = -80.23077000000001

$ws.Range("H97").Value = 1166.3334
$ws.Range("J97").Value = 749.5
$ws.Range("L97").Value = 2248.5
$ws.Range("N97").Value = -3240.5

$ws.Range("H116").Value = 222520.2
$ws.Range("I116").Value = 418930.12
$ws.Range("K116").Value = 418930.12
$ws.Range("M116").Value = -415488.12

$ws.Range("H118").Value = 709.3333
$ws.Range("I118").Value = 397.16666
$ws.Range("J118").Value = 1021.5
$ws.Range("K118").Value = 1191.49998
$ws.Range("L118").Value = 3064.5
$ws.Range("M118").Value = 465.5000199999999
$ws.Range("N118").Value = -6378.5

$ws.Range("H125").Value = 1633
$ws.Range("I125").Value = 1500
$ws.Range("J125").Value = 1699.5
$ws.Range("K125").Value = 13500
$ws.Range("L125").Value = 15295.5
$ws.Range("M125").Value = -11040
$ws.Range("N125").Value = -20215.5

$ws.Range("H137").Value = 4217.8945
$ws.Range("I137").Value = 3088.875
$ws.Range("J137").Value = 5039
$ws.Range("K137").Value = 9266.625
$ws.Range("L137").Value = 15117
$ws.Range("M137").Value = -6716.625
$ws.Range("N137").Value = -20217

$ws.Range("H138").Value = 2160.99
$ws.Range("I138").Value = 835.76
$ws.Range("J138").Value = 2602.7334
$ws.Range("K138").Value = 2507.28
$ws.Range("L138").Value = 7808.2002
$ws.Range("M138").Value = 2632.72
$ws.Range("N138").Value = -18088.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1071.2778
$ws.Range("I45").Value = 1056.8462
$ws.Range("J45").Value = 1108.8
$ws.Range("K45").Value = 1056.8462
$ws.Range("L45").Value = 1108.8
$ws.Range("M45").Value = -679.8462
$ws.Range("N45").Value = -1862.8

$ws.Range("H74").Value = 5001.524
$ws.Range("I74").Value = 4580.6313
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 4580.6313
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -3706.6313
$ws.Range("N74").Value = -10748

$ws.Range("H77").Value = 5001.524
$ws.Range("I77").Value = 4580.6313
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 22903.1565
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -18535.1565
$ws.Range("N77").Value = -53736

$ws.Range("H132").Value = 2032.3846
$ws.Range("I132").Value = 1128.6364
$ws.Range("J132").Value = 7003
$ws.Range("K132").Value = 3385.9092
$ws.Range("L132").Value = 21009
$ws.Range("M132").Value = -855.9092000000001
$ws.Range("N132").Value = -26069

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2205.1724
$ws.Range("I134").Value = 1748.2142
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 5244.642599999999
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -2709.642599999999
$ws.Range("N134").Value = -50070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 907.4167
$ws.Range("I22").Value = 299.8
$ws.Range("J22").Value = 1341.4286
$ws.Range("K22").Value = 299.8
$ws.Range("L22").Value = 1341.4286
$ws.Range("M22").Value = 50.19999999999999
$ws.Range("N22").Value = -2041.4286

$ws.Range("H31").Value = 17860344
$ws.Range("I31").Value = 1658.6666
$ws.Range("K31").Value = 1658.6666
$ws.Range("M31").Value = -1363.6666

$ws.Range("H34").Value = 17860344
$ws.Range("I34").Value = 1658.6666
$ws.Range("K34").Value = 1658.6666
$ws.Range("M34").Value = -1456.6666

$ws.Range("H122").Value = 1902.3334
$ws.Range("I122").Value = 966.7857
$ws.Range("K122").Value = 2900.3571
$ws.Range("M122").Value = -450.3571000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6667549.5
$ws.Range("I131").Value = 100000420
$ws.Range("J131").Value = 916.47144
$ws.Range("K131").Value = 300001260
$ws.Range("L131").Value = 2749.41432
$ws.Range("M131").Value = -299996220
$ws.Range("N131").Value = -12829.41432

$ws.Range("H140").Value = 2313.276
$ws.Range("I140").Value = 2394.524
$ws.Range("J140").Value = 2100
$ws.Range("K140").Value = 7183.572
$ws.Range("L140").Value = 6300
$ws.Range("M140").Value = -2003.572
$ws.Range("N140").Value = -16660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1106.6666
$ws.Range("I97").Value = 1010
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 1010
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -514
$ws.Range("N97").Value = -2292

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 34909.1
$ws.Range("I22").Value = 63685.312
$ws.Range("J22").Value = 2022
$ws.Range("K22").Value = 63685.312
$ws.Range("L22").Value = 2022
$ws.Range("M22").Value = -63390.312
$ws.Range("N22").Value = -2612

$ws.Range("H27").Value = 34909.1
$ws.Range("I27").Value = 63685.312
$ws.Range("J27").Value = 2022
$ws.Range("K27").Value = 63685.312
$ws.Range("L27").Value = 2022
$ws.Range("M27").Value = -63578.312
$ws.Range("N27").Value = -2236

$ws.Range("H46").Value = 3265.2
$ws.Range("I46").Value = 3380
$ws.Range("K46").Value = 3380
$ws.Range("M46").Value = -3192

$ws.Range("H93").Value = 4835069.5
$ws.Range("I93").Value = 8551016
$ws.Range("J93").Value = 4339
$ws.Range("K93").Value = 8551016
$ws.Range("L93").Value = 4339
$ws.Range("M93").Value = -8549768
$ws.Range("N93").Value = -6835

$ws.Range("H132").Value = 10643.849
$ws.Range("I132").Value = 9214.9
$ws.Range("J132").Value = 24933.334
$ws.Range("K132").Value = 27644.7
$ws.Range("L132").Value = 74800.00199999999
$ws.Range("M132").Value = -25114.7
$ws.Range("N132").Value = -79860.00199999999

$ws.Range("H136").Value = 4622.4443
$ws.Range("I136").Value = 1782.1818
$ws.Range("J136").Value = 9085.714
$ws.Range("K136").Value = 5346.5454
$ws.Range("L136").Value = 27257.142
$ws.Range("M136").Value = -2796.5454
$ws.Range("N136").Value = -32357.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2019
$ws.Range("I136").Value = 1354.4445
$ws.Range("K136").Value = 4063.3335
$ws.Range("M136").Value = -1513.3335
